$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; A="Intestazione: Settore, Servizio e Autorità che emana il provvedimento e tipo di provvedimento. Nella determina oltre il settore va indicato il servizio di appartenenza presente nella macrostruttura vigente."; B="RISPOSTA: NO`nNote: Non è esplicitato nella determina il servizio di appartenenza presente nella macrostruttura vigente."; C="NO"},
    @{Row=3; A="Oggetto della determinazione"; B="RISPOSTA : NO`nNote: L’oggetto della determinazione menziona il codice 15275 anziché 15271, è necessario correggerlo."; C="NO"},
    @{Row=4; A="Indicazione del CIG e /o del CUP del CUI, del CPV"; B="RISPOSTA: SI"; C="SI"},
    @{Row=5; A="Centro di Costo: si indica il centro di costo al quale imputare la spesa."; B="RISPOSTA : SI"; C="SI"},
    @{Row=6; A="Estremi decreto sindacale di nomina del dirigente."; B="RISPOSTA : SI"; C="SI"},
    @{Row=7; A="Estremi della delega alla firma dell'atto, se persona diversa dal Dirigente."; B="RISPOSTA: NO"; C="NO"},
    @{Row=8; A="Estremi atto di nomina del Responsabile del Progetto."; B="RISPOSTA : SI"; C="SI"},
    @{Row=9; A="Conflitto d'interessi "; B="RISPOSTA: SI`nNota: La determina menziona che non sussistono situazioni di conflitto di interessi in relazione ai soggetti coinvolti nella procedura, in conformità all'art. 16 del D.lgs. n. 36/2023."; C="SI"},
    @{Row=10; A="**Normativa specifica**:`nCodice dei contratti pubblici (D.Lgs. n. 36.2023.);`nL.R. n. 8/2018 e ss.mm.ii. (per quanto applicabile);`nArt. 1, D.L. 6 luglio 2012 n. 95 convertito dalla legge 7 agosto 2012, n. 135 (adesione convenzione Consip);`nArt. 1, comma 450, della legge 27 Dicembre 2006, n. 296 e ss.mm.ii. (per quanto applicabile);`nArt. 26, della legge n. 488/1999 e ss.mm.ii.;`nDPCM 24 dicembre 2015."; B="RISPOSTA : SI"; C="SI"},
    @{Row=11; A="**Normativa generale**:`nTUEL;`nLegge n. 241/90;`nDPR n. 62/2013 Codice comportamento dipendenti pubblici;`nL. n. 190/2012;`nD.Lgs. n. 33/2013;`nLegge n.136/2010."; B="RISPOSTA : SI"; C="SI"},
    @{Row=12; A="**Norme di principio**`nArt. 4 del Dlgs n.36/2023. (Criterio interpretativo e applicativo)"; B="RISPOSTA : NO"; C="NO"},
    @{Row=13; A="Regolamenti dell’ente quali:`n- Statuto Comunale;`n- Regolamento di contabilità;`n- Regolamento dei contratti"; B="RISPOSTA: NO`nNote: Nella determina non si fa esplicito riferimento allo Statuto Comunale, al Regolamento di contabilità, né al Regolamento dei contratti."; C="NO"},
    @{Row=14; A="Termini per la conclusione della procedura "; B="RISPOSTA: NO`nNote: Non è presente alcun riferimento ai termini per la conclusione della procedura nella determina."; C="NO"},
    @{Row=15; A="Fine che con il contratto si intende perseguire"; B="RISPOSTA: SI"; C="SI"},
    @{Row=16; A="Oggetto del contratto"; B="RISPOSTA: NO `nNote: La determina non specifica un capitolato d’appalto, quaderno d’oneri o schema di contratto."; C="NO"},
    @{Row=17; A="DUVRI (per i servizi)"; B="RISPOSTA: SI"; C="SI"},
    @{Row=18; A="Costo della mano d’opera; contratto applicabile "; B="NO"; C="NO"},
    @{Row=19; A="Suddivisione in lotti"; B="RISPOSTA : NO`nNota: La determina non menziona la motivazione sulla mancata suddivisione dell'appalto in lotti ai sensi dell’art. 58 del Dlgs n. 36/2023."; C="NO"},
    @{Row=20; A="Validazione (in caso di lavori pubblici)"; B="RISPOSTA : NO`nNote: Nel testo della determina non è menzionata la validazione del progetto posto a base di gara."; C="NO"},
    @{Row=21; A="Criteri Ambientali minimi"; B="RISPOSTA: NO`nNote: Nella determina non è presente alcun riferimento ai Criteri Ambientali Minimi (CAM) o alla loro applicazione."; C="NO"},
    @{Row=22; A="Forma del contratto"; B="RISPOSTA : SI"; C="SI"},
    @{Row=23; A="Clausole del contratto ritenute essenziali"; B="RISPOSTA: SI"; C="SI"},
    @{Row=24; A="Cauzione provvisoria"; B="RISPOSTA: NO`nNote: Non viene richiesta la garanzia definitiva in quanto l’affidamento non lo prevede, secondo quanto indicato nella determina."; C="NO"},
    @{Row=25; A="Modalità di scelta del contraente e ragioni che ne sono alla base (motivare in modo specifico) e criteri di selezione delle offerte"; B="RISPOSTA: NO`nNote: La determina non fornisce una motivazione specifica riguardante la modalità di scelta del contraente e i criteri di selezione delle offerte, come richiesto dal punto 24."; C="NO"},
    @{Row=26; A="Rotazione degli affidamenti "; B="RISPOSTA : NO`nNote: La determina specifica che il principio di rotazione di cui all’art. 49 del d.lgs. 36/2023 può essere derogato per gli affidamenti diretti di importo inferiore a 5.000 euro. Questo implica che non si applica la rotazione per la presente procedura."; C="NO"},
    @{Row=27; A="Indicazione dell’importo massimo stimato a base di gara`nQuadro economico dell’intervento  comprensivo di `nImporto a base di gara, oneri fiscali, costi della sicurezza/interferenziali, revisione dei prezzi, somme a disposizione, supporto al RUP, incentivi (art. 45) modificazioni del contratto preventivabili, proroga, lavori o servizi analoghi, quinto d’obbligo"; B="RISPOSTA: NO`nNota: Nella determina non è presente un'indicazione chiara dell'importo massimo stimato a base di gara e del quadro economico dell'intervento, come i costi della sicurezza, oneri fiscali, e altre voci richieste. Questo rappresenta una grave omissione."; C="NO"},
    @{Row=28; A="Riferimenti all’obbligo di utilizzo degli strumenti di acquisto e di negoziazione messi a disposizione da CONSIP e da soggetti aggregatori "; B="RISPOSTA : NO`nNote: La determina non menziona chiaramente il rispetto degli obblighi di utilizzo degli strumenti di acquisto e di negoziazione messi a disposizione da CONSIP e da soggetti aggregatori, come richiesto."; C="NO"},
    @{Row=29; A="Riferimenti alla qualificazione del Comune quale Stazione appaltante"; B="RISPOSTA : NO"; C="NO"},
    @{Row=30; A="Riferimenti alla deliberazione di approvazione del bilancio."; B="RISPOSTA: NO`nNote: Non sono presenti riferimenti alla deliberazione di approvazione del bilancio pluriennale vigente nella determina."; C="NO"},
    @{Row=31; A="Riferimenti deliberazione approvazione PEG e del PIAO"; B="RISPOSTA : NO"; C="NO"},
    @{Row=32; A="Impegno di spesa/prenotazione/copertura finanziaria"; B="RISPOSTA : SI `nNota: L'importo da impegnare è di € 4.026,00, assoggettato a IVA al 22% per una parte (aggiornamento Software Auditing PA) mentre la formazione è esente da IVA ai sensi dell’articolo 10, primo comma, n. 20) del D.P.R. 633/72. Non ci sono riferimenti a bilancio e PEG non approvati o impegni di spesa pluriennali nel testo della determina."; C="SI"},
    @{Row=33; A="Conformità dei pagamenti con le regole della finanza pubblica vigenti."; B="RISPOSTA : SI"; C="SI"},
    @{Row=34; A="Contributo ANAC"; B="RISPOSTA : NO"; C="NO"},
    @{Row=35; A="Richiamati e/o allegati"; B="RISPOSTA : SI"; C="SI"},
    @{Row=36; A="Obblighi di pubblicità e trasparenza"; B="RISPOSTA: SI"; C="SI"},
    @{Row=37; A="Pubblicazione atti all’Albo Pretorio Online"; B="RISPOSTA: SI"; C="SI"},
    @{Row=38; A="Acquisizione del visto di regolarità contabile ai sensi dell’art. 153 del D.Lgs. n. 267/2000"; B="RISPOSTA : SI"; C="SI"},
    @{Row=39; A="Sottoscrizione del Dirigente firmatario dell’atto e dove previsto del R.P. o RUP."; B="RISPOSTA: SI"; C="SI"},
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 1).Value = $item.A
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
}

# Remove rows 40-47 (no longer present in the updated checklist)
$ws.Range("A40:C47").EntireRow.Delete()

